$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.998.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.34%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.031.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.02%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'594.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.14%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'153.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +8.02%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.028.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.32%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'6.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +17.64%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.50%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +3.08%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +3.96%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'35.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.71%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  -0.24%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.539.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.34%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +3.73%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'62.918.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.09%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'3.029.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.16%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'452.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.79%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'14.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.40%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +2.97%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +3.88%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'83.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.46%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'11.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +11.25%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +8.19%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'12.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.77%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -0.07%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +13.22%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'7.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.60%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +1.81%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +0.07%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'27.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.58%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +3.67%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0₃0864"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.74%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +3.07%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +2.81%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +11.55%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +9.12%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.01%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'50.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.04%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +1.70%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +16.95%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'44.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +15.14%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'392.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.98%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +4.07%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'2.720.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.17%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'132.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.59%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'25.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +10.72%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -0.01%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +8.28%  "
$ws.Range("E51").Style = "Normal"
